$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "91.623.51"
Set-TextValue "E2" "  +1.08%  "

Set-TextValue "D3" "3.134.57"
Set-TextValue "E3" "  -0.03%  "

Set-TextValue "D4" "0.999"
Set-TextValue "E4" "  -0.09%  "

Set-TextValue "D5" "242.33"
Set-TextValue "E5" "  +1.47%  "

Set-TextValue "D6" "628.36"
Set-TextValue "E6" "  -0.63%  "

Set-TextValue "E7" "  +9.51%  "

Set-TextValue "D8" "0.377"
Set-TextValue "E8" "  +6.00%  "

Set-TextValue "D9" "0.999"
Set-TextValue "E9" "  -0.13%  "

Set-TextValue "D10" "3.135.21"
Set-TextValue "E10" "  +0.09%  "

Set-TextValue "D11" "0.770"
Set-TextValue "E11" "  +6.71%  "

Set-TextValue "E12" "  +4.61%  "

Set-TextValue "D13" "0.0000254"
Set-TextValue "E13" "  +4.64%  "

Set-TextValue "D14" "36.02"
Set-TextValue "E14" "  -1.53%  "

Set-TextValue "D15" "5.53"
Set-TextValue "E15" "  -1.62%  "

Set-TextValue "D16" "91.131.45"
Set-TextValue "E16" "  +0.68%  "

Set-TextValue "D17" "3.725.27"
Set-TextValue "E17" "  +0.42%  "

Set-TextValue "D18" "3.158.12"
Set-TextValue "E18" "  +1.16%  "

Set-TextValue "D19" "3.79"
Set-TextValue "E19" "  +3.38%  "

Set-TextValue "D20" "14.72"
Set-TextValue "E20" "  +2.18%  "

Set-TextValue "D21" "0.0000215"
Set-TextValue "E21" "  +1.51%  "

Set-TextValue "D22" "5.88"
Set-TextValue "E22" "  +3.08%  "

Set-TextValue "D23" "454.05"
Set-TextValue "E23" "  +0.34%  "

Set-TextValue "D24" "9.18"
Set-TextValue "E24" "  +0.90%  "

Set-TextValue "D25" "5.98"
Set-TextValue "E25" "  +2.85%  "

Set-TextValue "D26" "93.55"
Set-TextValue "E26" "  +3.29%  "

Set-TextValue "D27" "12.04"
Set-TextValue "E27" "  -3.47%  "

Set-TextValue "D28" "3.297.83"

Set-TextValue "E29" "  +0.10%  "

Set-TextValue "D30" "0.180"
Set-TextValue "E30" "  +11.30%  "

Set-TextValue "D31" "0.124"
Set-TextValue "E31" "  +42.41%  "

Set-TextValue "E32" "  +14.84%  "

Set-TextValue "D33" "9.14"
Set-TextValue "E33" "  -9.44%  "

Set-TextValue "D34" "1.01"
Set-TextValue "E34" "  +26.15%  "

Set-TextValue "E35" "  +9.02%  "

Set-TextValue "D36" "27.04"
Set-TextValue "E36" "  -1.91%  "

Set-TextValue "D37" "7.68"
Set-TextValue "E37" "  +7.35%  "

Set-TextValue "D38" "4.22"
Set-TextValue "E38" "  +25.26%  "

Set-TextValue "D39" "503.59"
Set-TextValue "E39" "  -2.20%  "

Set-TextValue "D40" "1.94"
Set-TextValue "E40" "  +0.06%  "

Set-TextValue "D41" "3.64"
Set-TextValue "E41" "  -6.17%  "

Set-TextValue "E42" "  -0.61%  "

Set-TextValue "D43" "0.427"
Set-TextValue "E43" "  -0.44%  "

Set-TextValue "E44" "  -0.09%  "

Set-TextValue "E45" "  +0.01%  "

Set-TextValue "D46" "1.95"
Set-TextValue "E46" "  -0.20%  "

Set-TextValue "D47" "157.28"
Set-TextValue "E47" "  +5.22%  "

Set-TextValue "E48" "  -0.65%  "

Set-TextValue "D49" "4.59"
Set-TextValue "E49" "  +0.44%  "

Set-TextValue "B50" "ImmutableX"
Set-TextValue "C50" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue "D50" "1.35"
Set-TextValue "E50" "  +0.22%  "

Set-TextValue "B51" "OKB"
Set-TextValue "C51" "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D51" "45.08"
Set-TextValue "E51" "  -0.92%  "
